# Auto-generated script to apply the cryptos.xlsx price/volume refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.893.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "'1.551.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.14%  "

$ws.Range("E4").Value = "  +0.54%  "

$ws.Range("D5").Value = "'206.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  +0.57%  "

$ws.Range("D8").Value = "'21.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.42%  "

$ws.Range("E9").Value = "  +1.27%  "

$ws.Range("E10").Value = "  +1.04%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "'1.773.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("D13").Value = "'1.551.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("D15").Value = "'0.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.64%  "

$ws.Range("D16").Value = "'26.896.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").Value = "'61.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("E18").Value = "  +2.00%  "

$ws.Range("D19").Value = "'0.0₃0688"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("D20").Value = "'7.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

$ws.Range("E22").Value = "  +0.81%  "

$ws.Range("E23").Value = "  +1.21%  "

$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").Value = "'153.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "

$ws.Range("D26").Value = "'6.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("D27").Value = "'14.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("E29").Value = "  +0.99%  "

$ws.Range("E30").Value = "  +2.69%  "

$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D33").Value = "'1.422.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.82%  "

$ws.Range("E34").Value = "  +2.90%  "

$ws.Range("E35").Value = "  +3.88%  "

$ws.Range("E36").Value = "  +2.03%  "

$ws.Range("E37").Value = "  +1.12%  "

$ws.Range("E38").Value = "  +0.86%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("D42").Value = "'5.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("E44").Value = "  +3.78%  "

$ws.Range("D45").Value = "'63.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").Value = "'1.686.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D48").Value = "'86.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("D49").Value = "'0.0522"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.07%  "

$ws.Range("E50").Value = "  +4.74%  "

$ws.Range("D51").Value = "'0.0954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "
